# Energy Monitoring ISO 50001 template update
# Restructures the "Engine" and "Air Volume" sheets from a wide (1 header row
# + 1 data row, 4 compressor columns) layout into a tall (2 columns, 5 rows)
# layout suitable for tile views.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Engine"
# ---------------------------------------------------------------------------
$wsEngine = $wb.Worksheets.Item("Engine")
$wsEngine.Cells.Clear()

$wsEngine.Range("A1").Value = "Compressor"

$wsEngine.Range("A2").Value = "Compressor 1"
$wsEngine.Range("B2").Value = "y"

$wsEngine.Range("A3").Value = "Compressor 2"
$wsEngine.Range("B3").Value = "n"

$wsEngine.Range("A4").Value = "Compressor 3"
$wsEngine.Range("B4").Value = "y"

$wsEngine.Range("A5").Value = "Compressor 4"
$wsEngine.Range("B5").Value = "n"

# ---------------------------------------------------------------------------
# Sheet 2: "Air Volume"
# ---------------------------------------------------------------------------
$wsVolume = $wb.Worksheets.Item("Air Volume")
$wsVolume.Cells.Clear()

$wsVolume.Range("A1").Value = "Compressor"
$wsVolume.Range("B1").Value = "Volume"

$wsVolume.Range("A2").Value = "Compressor 1"
$wsVolume.Range("B2").Value = 777

$wsVolume.Range("A3").Value = "Compressor 2"
$wsVolume.Range("B3").Value = 0

$wsVolume.Range("A4").Value = "Compressor 3"
$wsVolume.Range("B4").Value = 2270

$wsVolume.Range("A5").Value = "Compressor 4"
$wsVolume.Range("B5").Value = 0

$wsVolume.Columns.Item(1).ColumnWidth = 14.830729166666666

$wsVolume.Range("A6:D17").Select()

# Now fill in the "Turned on" header on the Engine sheet, after "Volume" has
# already been registered as a shared string, to match the shared string
# table ordering of the saved workbook.
$wsEngine.Range("B1").Value = "Turned on"
$wsEngine.Columns.Item(1).ColumnWidth = 12.053385416666666
$wsEngine.Range("B1").Select()

# ---------------------------------------------------------------------------
# Restore the originally active sheet ("Energy Consumption") so the workbook
# still opens on the same tab it did before the edit.
# ---------------------------------------------------------------------------
$wsEnergy = $wb.Worksheets.Item("Energy Consumption")
$wsEnergy.Activate()
